# Feedback form basic layout
# Adds a new "feedback" service row (row 38) to the Services list sheet,
# mirroring the structure/formulas of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (37) down to the new row (38)
$ws.Range("B37:O37").Copy() | Out-Null
$ws.Range("B38:O38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Values for the new row ----
# NOTE: assignment order chosen to reproduce the original authoring order of
# the shared-string table (feedback, WS-FED-01, app.feedback.question.get,
# Get feedback Question, /getquestion).
$ws.Cells.Item(38, 2).Value2 = "feedback"                       # B38 Application
$ws.Cells.Item(38, 4).Value2 = "WS-FED-01"                      # D38 ServiceCode
$ws.Cells.Item(38, 5).Value2 = "app.feedback.question.get"      # E38 queryId
$ws.Cells.Item(38, 3).Value2 = "Get feedback Question"          # C38 Service Name
$ws.Cells.Item(38, 6).Value = "'false"                          # F38 logActivity (text "false")
$ws.Cells.Item(38, 7).Value2 = "feedback"                        # G38 BasePath
$ws.Cells.Item(38, 8).Value2 = "/getquestion"                   # H38 servicePath
$ws.Cells.Item(38, 9).Value2 = "POST"                            # I38 ServiceType

# F38 picked up a "quoted text" style variant from the apostrophe-prefixed
# assignment above; re-apply the plain text-format styling from F37 so the
# cell style matches the rest of the column.
$ws.Range("F37").Copy() | Out-Null
$ws.Range("F38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Formulas for the new row (mirrors rows above) ----
$ws.Cells.Item(38, 13).Formula = '=CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D38,"'',''CONNON_CONFIG'', 0, ''",C38,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Cells.Item(38, 14).Formula = '=CONCAT(IF(I38="GET","@GetMapping(",IF(I38="POST","@PostMapping(",IF(I38="DELETE","@DeleteMapping(",IF(I38="PUT","@PutMapping(","")))),CHAR(34),H38,CHAR(34),")")'
$ws.Cells.Item(38, 15).Formula = '=CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D38,,CHAR(34),", serviceName = ",CHAR(34),C38,CHAR(34), ", queryId = ",CHAR(34),E38,CHAR(34),", logActivity =",F38,")")'

# ---- Update the view: scroll position and active selection ----
# (the existing frozen pane at row 3 is left untouched; we just move the
# viewport/selection the way the author did while reviewing the new row)
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 21
$ws.Range("H43").Select() | Out-Null

$wb.Save()
